$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56

$ws.Range("A$row`:D$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-07"
$ws.Cells.Item($row, 2).Value = "09:17:58"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "05"

$ws.Range("A$row`:D$row").ClearFormats()

$ws.Cells.Item($row, 5).Value = 125826
$ws.Cells.Item($row, 6).Value = 141760
$ws.Cells.Item($row, 7).Value = 167556
$ws.Cells.Item($row, 8).Value = 158043
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142843
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191336
$ws.Cells.Item($row, 14).Value = 115171
$ws.Cells.Item($row, 15).Value = 44680
$ws.Cells.Item($row, 16).Value = 28272
$ws.Cells.Item($row, 17).Value = 63397
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 40255
$ws.Cells.Item($row, 20).Value = -1
